$d = $word.ActiveDocument

$replacements = @(
    @{old="322×6=1932"; new="957×6=5742"},
    @{old="213×6=1278"; new="555×7=3885"},
    @{old="353×2=706";  new="355×4=1420"},
    @{old="359×6=2154"; new="638×6=3828"},
    @{old="572×5=2860"; new="129×9=1161"},
    @{old="170×6=1020"; new="399×7=2793"},
    @{old="475×2=950";  new="523×6=3138"},
    @{old="251×7=1757"; new="686×4=2744"},
    @{old="849×2=1698"; new="773×8=6184"},
    @{old="963×6=5778"; new="657×7=4599"},
    @{old="286×9=2574"; new="369×8=2952"},
    @{old="914×2=1828"; new="440×5=2200"},
    @{old="805×5=4025"; new="204×4=816"},
    @{old="780×3=2340"; new="180×7=1260"},
    @{old="540×3=1620"; new="752×5=3760"},
    @{old="512×3=1536"; new="889×5=4445"},
    @{old="785×9=7065"; new="893×5=4465"},
    @{old="842×4=3368"; new="684×3=2052"},
    @{old="260×5=1300"; new="952×7=6664"},
    @{old="856×3=2568"; new="242×3=726"},
    @{old="857×3=2571"; new="419×8=3352"},
    @{old="583×5=2915"; new="949×4=3796"},
    @{old="759×4=3036"; new="178×6=1068"},
    @{old="249×5=1245"; new="522×5=2610"},
    @{old="630×6=3780"; new="262×4=1048"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
